$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# STEP A/B: split the "Alexandre et Guillaume..." paragraph into two:
#   - a shorter "Alexandre..." paragraph ending at "...sur Firestore."
#   - a new "et Guillaume..." paragraph about the scrolling background
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Alexandre et Guillaume s’occupent du jeu avec le barbecue et de faire un objet score à sortir pour pouvoir enregistrer le meilleur score sur Firestore + faire un classement local de ses propres scores + DAO des données en local (avec localstorage).")
$r.Text = "Alexandre s’occupe du jeu avec le barbecue et de faire un objet score à sortir pour pouvoir enregistrer le meilleur score sur Firestore."
$r.InsertParagraphAfter()
$newPara = $r.Next(4, 1)

$r2 = $d.Content
$r2.Find.Execute("Alexandre s’occupe du jeu avec le barbecue et de faire un objet score à sortir pour pouvoir enregistrer le meilleur score sur Firestore.")
$guillaumeRange = $r2.Paragraphs(1).Next().Range
$guillaumeRange.Text = "et Guillaume s’occupe du fond défilant derrière le menu ainsi que de la page de chargement + aide avec Alexandre."

# ---------------------------------------------------------------------
# STEP C: rewrite the "et fond défilant derrière le menu." paragraph
# (now located right after "Simon s’occupe...") into the new
# "S’il y a le temps..." paragraph about the local ranking / DAO.
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("et fond défilant derrière le menu.")
$r3.Text = "S’il y a le temps faire un classement local de ses propres scores + DAO des données en local (avec localstorage)."

# ---------------------------------------------------------------------
# STEP E: add "-case pour mettre le pseudo" after the "Jouer" bullet
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("-bouton Jouer pour lancer la partie")
$r4.InsertParagraphAfter()
$r5 = $d.Content
$r5.Find.Execute("-bouton Jouer pour lancer la partie")
$pseudoRange = $r5.Paragraphs(1).Next().Range
$pseudoRange.Text = "-case pour mettre le pseudo"

# ---------------------------------------------------------------------
# STEP F: rewrite "Page de classement local :" -> "Page de classement
# local (peut-être):" and add the new bullet under it.
# ---------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.Execute("Page de classement local :")
$r6.Text = "Page de classement local (peut-être):"
$r6.InsertParagraphAfter()
$r7 = $d.Content
$r7.Find.Execute("Page de classement local (peut-être):")
$classementRange = $r7.Paragraphs(1).Next().Range
$classementRange.Text = "-classement de tous les scores obtenus sur l’appareil"

# ---------------------------------------------------------------------
# STEP G: add "- Liste des meilleurs scores..." after
# "Page de classement Firestore :" (end of document).
# ---------------------------------------------------------------------
$r8 = $d.Content
$r8.Find.Execute("Page de classement Firestore :")
$r8.InsertParagraphAfter()
$r9 = $d.Content
$r9.Find.Execute("Page de classement Firestore :")
$firestoreRange = $r9.Paragraphs(1).Next().Range
$firestoreRange.Text = "- Liste des meilleurs scores obtenus pour chaque joueur"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
